# Instructions.pptx edit: "Finished Alchemy practice 1 and 2"
#
# 1) Remove the last two slides (Alchemy practice 1 & 2) from the deck.
# 2) On slide 4, tidy up the wording of the last bullet in "TextBox 6"
#    (split it into 3 runs, fixing a "seperate" typo in the process) and
#    shrink the text box now that there is less text in it.

$p = $ppt.ActivePresentation

# --- 1) Drop the trailing two slides -----------------------------------
while ($p.Slides.Count -gt 5) {
    $p.Slides.Item($p.Slides.Count).Delete()
}

# --- 2) Slide 4 bullet rewrite + resize ---------------------------------
$slide4 = $p.Slides.Item(4)
$box = $slide4.Shapes.Item("TextBox 6")

$tr = $box.TextFrame.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)

# Step 1: land the new sentence in one shot. The host diff-patches this
# against the old run, so the common "Each experiment is a " prefix and
# the (non-existent, this time) common suffix stay their own runs and
# only the changed middle becomes a new run.
$lastPara.Text = "Each experiment is a distinct alchemical set, seperate from other experiments"

# Step 2: re-stamp characters 1-46 ("Each experiment is a distinct " +
# "alchemical set, ") with identical text -- this merges the prefix run
# together with the leading part of the middle run into one clean run
# without altering any characters.
$lastPara.Characters(1, 46).Text = "Each experiment is a distinct alchemical set, "

# Step 3: likewise re-stamp characters 55-77 (" from other experiments")
# so the tail merges into a single trailing run, leaving exactly
# "seperate" isolated as its own (8-character) middle run.
$lastPara.Characters(55, 23).Text = " from other experiments"

# Shrink the box now that the paragraph is shorter (off/cx unchanged).
$box.Height = 298.08284
